# Update the "Content of the dataset" description on the Metadata sheet
# to reflect that candidate periods are now defined relative to the
# DU pregnant persontime rather than the (older) "pregnancy cohort" wording.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "This applies some of the steps in the appendix: selects all the periods when women with MS are outside of the DU pregnant persontime "
